# Add 7 new data rows (16-22) to the Landscaping Data sheet, matching the
# rows already present (same collection date 45789 / 2025-05-12, same
# weather figures), and update the sheet's view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout (row 1 headers):
# A Date | B Plant_Type | C Plant_Size | D Low | E High | F Temp_Diff (=E-D)
# G Rain | H Growth | I Pruned | J Quadrant | K Shade | L UV | M Humidity
# N Dew_Point | O Pressure | P Wind_Gust | Q Cloud_Cover | R Visibility
# S AQI | T Pollen

$newRows = @(
    @{ Date=45789; Plant_Type="Flowering";    Plant_Size="Large";  Low=64; High=80; Rain=0.33; Growth=1;   Pruned="No"; Quadrant=2; Shade="Neutral"; UV=8; Humidity=0.42; Dew_Point=51; Pressure=30.2; Wind_Gust=18; Cloud_Cover=0.36; Visibility=9.9; AQI=53; Pollen=45 },
    @{ Date=45789; Plant_Type="Nonflowering"; Plant_Size="Medium"; Low=64; High=80; Rain=0.33; Growth=0.2; Pruned="No"; Quadrant=3; Shade="Dark";    UV=8; Humidity=0.42; Dew_Point=51; Pressure=30.2; Wind_Gust=18; Cloud_Cover=0.36; Visibility=9.9; AQI=53; Pollen=45 },
    @{ Date=45789; Plant_Type="Nonflowering"; Plant_Size="Small";  Low=64; High=80; Rain=0.33; Growth=0.3; Pruned="No"; Quadrant=3; Shade="Neutral"; UV=8; Humidity=0.42; Dew_Point=51; Pressure=30.2; Wind_Gust=18; Cloud_Cover=0.36; Visibility=9.9; AQI=53; Pollen=45 },
    @{ Date=45789; Plant_Type="Nonflowering"; Plant_Size="Medium"; Low=64; High=80; Rain=0.33; Growth=0.5; Pruned="No"; Quadrant=3; Shade="Bright";  UV=8; Humidity=0.42; Dew_Point=51; Pressure=30.2; Wind_Gust=18; Cloud_Cover=0.36; Visibility=9.9; AQI=53; Pollen=45 },
    @{ Date=45789; Plant_Type="Nonflowering"; Plant_Size="Medium"; Low=64; High=80; Rain=0.33; Growth=0.5; Pruned="No"; Quadrant=3; Shade="Bright";  UV=8; Humidity=0.42; Dew_Point=51; Pressure=30.2; Wind_Gust=18; Cloud_Cover=0.36; Visibility=9.9; AQI=53; Pollen=45 },
    @{ Date=45789; Plant_Type="Nonflowering"; Plant_Size="Large";  Low=64; High=80; Rain=0.33; Growth=0.1; Pruned="No"; Quadrant=4; Shade="Dark";    UV=8; Humidity=0.42; Dew_Point=51; Pressure=30.2; Wind_Gust=18; Cloud_Cover=0.36; Visibility=9.9; AQI=53; Pollen=45 },
    @{ Date=45789; Plant_Type="Tree";         Plant_Size="Medium"; Low=64; High=80; Rain=0.33; Growth=1.5; Pruned="No"; Quadrant=1; Shade="Neutral"; UV=8; Humidity=0.42; Dew_Point=51; Pressure=30.2; Wind_Gust=18; Cloud_Cover=0.36; Visibility=9.9; AQI=53; Pollen=45 }
)

$startRow = 16
$lastDataRow = 15

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Copy the last existing row's per-cell formatting (keeps column A's
    # short-date number format, and everything else unstyled) down onto the
    # new row before filling in values.
    $ws.Range("A$lastDataRow`:T$lastDataRow").Copy($ws.Range("A$r`:T$r"))

    $ws.Cells.Item($r, 1).Value = $row.Date
    $ws.Cells.Item($r, 2).Value = $row.Plant_Type
    $ws.Cells.Item($r, 3).Value = $row.Plant_Size
    $ws.Cells.Item($r, 4).Value = $row.Low
    $ws.Cells.Item($r, 5).Value = $row.High
    $ws.Cells.Item($r, 6).Formula = "=E$r-D$r"
    $ws.Cells.Item($r, 7).Value = $row.Rain
    $ws.Cells.Item($r, 8).Value = $row.Growth
    $ws.Cells.Item($r, 9).Value = $row.Pruned
    $ws.Cells.Item($r, 10).Value = $row.Quadrant
    $ws.Cells.Item($r, 11).Value = $row.Shade
    $ws.Cells.Item($r, 12).Value = $row.UV
    $ws.Cells.Item($r, 13).Value = $row.Humidity
    $ws.Cells.Item($r, 14).Value = $row.Dew_Point
    $ws.Cells.Item($r, 15).Value = $row.Pressure
    $ws.Cells.Item($r, 16).Value = $row.Wind_Gust
    $ws.Cells.Item($r, 17).Value = $row.Cloud_Cover
    $ws.Cells.Item($r, 18).Value = $row.Visibility
    $ws.Cells.Item($r, 19).Value = $row.AQI
    $ws.Cells.Item($r, 20).Value = $row.Pollen
}

# Match the saved selection/active cell from the edited workbook.
[void]$ws.Range("Q23").Select()
